$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously-filled step-2 comment in row 16 (content moved into the
# newly added test cases below; the row now just closes out test case 4).
$ws.Range("H16").Value = ""

# --- Test case 5: CRUD goals ---
$ws.Range("A18").Value = 5
$ws.Range("B18").Value = 'CRUD goals'
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 'Try to create, read, update and delete a goal'
$ws.Range("E18").Value = 'Can do CRUD operations on goals with ease'
$ws.Range("F18").Value = 'Error: create is missing data for required field'
$ws.Range("G18").Value = 'Failed'
$ws.Range("H18").Value = 'No default created_at time. This is causing the error.'

$ws.Range("F19").Value = 'Error: null value in column "goal_type" violates not-null constraint'
$ws.Range("G19").Value = 'Failed'
$ws.Range("H19").Value = 'Goal type is missing'

$ws.Range("C20").Value = 2
$ws.Range("F20").Value = 'AttributeError: type object ''Goal'' has no attribute ''week_start'''
$ws.Range("G20").Value = 'Failed'
$ws.Range("H20").Value = 'I need to change week_start to created'

$ws.Range("F21").Value = 'Weekly expense source is now easily manipulated by owner'
$ws.Range("G21").Value = 'Passed'

$ws.Range("H22").Value = 'Also email must be set to unique, otherwise multiple users may be returned who have the same email. '

$ws.Range("A23").Value = 6
$ws.Range("B23").Value = 'CRUD exercise log item'
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 'Try to create, read, update and delete exercise log items'
$ws.Range("E23").Value = 'Can do CRUD operations on exercise log item'
$ws.Range("F23").Value = 'Error: exercise log item has no property ''username'''
$ws.Range("G23").Value = 'Failed'
$ws.Range("H23").Value = 'Forgot to set the username field'

$ws.Range("F24").Value = 'invalid input syntax for type integer: "test2@gmail.com"'
$ws.Range("G24").Value = 'Failed'

$ws.Range("F25").Value = 'AttributeError: type object ''ExerciseLogItem'' has no attribute ''created'''
$ws.Range("G25").Value = 'Failed'

$ws.Range("G26").Value = 'Passed'

$ws.Range("H27").Value = 'Also email must be set to unique, otherwise multiple users may be returned who have the same email. '

$ws.Range("A28").Value = 7
$ws.Range("B28").Value = 'Images'
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 'Try to upload new images, see them in app and delete them'
$ws.Range("E28").Value = 'Can upload, view and delete images'
$ws.Range("F28").Value = '(psycopg2.errors.NotNullViolation) null value in column "user_id" violates not-null constraint'
$ws.Range("G28").Value = 'Failed'
$ws.Range("H28").Value = 'have not set user_id field'

$ws.Range("F29").Value = 'TypeError: user_image_show() got an unexpected keyword argument ''id'''
$ws.Range("G29").Value = 'Failed'
$ws.Range("H29").Value = 'argument names in the image end point and the function must be consistent. They are image_number and id'

$ws.Range("F30").Value = 'AttributeError: ''Image'' object has no attribute ''username'''
$ws.Range("G30").Value = 'Failed'
$ws.Range("H30").Value = 'need to change username to user_id'

$ws.Range("G31").Value = 'Passed'

$ws.Range("H32").Value = 'Also email must be set to unique, otherwise multiple users may be returned who have the same email. '

$ws.Range("A33").Value = 8
$ws.Range("B33").Value = 'Mental Health Surveys'
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 'Try CRUD operations on surveys'
$ws.Range("E33").Value = 'Can create, read, update and delete surveys'
$ws.Range("F33").Value = 'update or delete on table "mental_health_surveys" violates foreign key constraint "survey_questions_mental_health_survey_id_fkey" on table "survey_questions"'
$ws.Range("G33").Value = 'Failed'
$ws.Range("H33").Value = 'need to check for any survey questions that are dependent on mental health surveys before deleting a mental health survey'

$ws.Range("G34").Value = 'Passed'

$ws.Range("H35").Value = 'Also email must be set to unique, otherwise multiple users may be returned who have the same email. '

$ws.Range("A36").Value = 9
$ws.Range("B36").Value = 'Questions'
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 'Try CRUD operations on the questions'
$ws.Range("E36").Value = 'Can create, read, update and delete questions'
$ws.Range("F36").Value = 'missing id and question_text'
$ws.Range("G36").Value = 'Failed'
$ws.Range("H36").Value = 'question schema and question model do not match. Need to get rid of id and replace "question_text" with just "text"'

$ws.Range("F37").Value = 'AttributeError: ''Question'' object has no attribute ''update'''
$ws.Range("G37").Value = 'Failed'
$ws.Range("H37").Value = 'don''t use the .first() method on the question object'

$ws.Range("G38").Value = 'Passed'

$ws.Range("H39").Value = 'Also email must be set to unique, otherwise multiple users may be returned who have the same email. '

$ws.Range("A40").Value = 10
$ws.Range("B40").Value = 'Questions'
$ws.Range("C40").Value = 1
$ws.Range("D40").Value = 'Try CRUD operations on the survey questions'
$ws.Range("E40").Value = 'Can create, read, update and delete survey questions'
$ws.Range("F40").Value = 'Cannot see the survey details and question details in the survey question object in json'
$ws.Range("G40").Value = 'Failed'
$ws.Range("H40").Value = 'Need to alter the schemas and models so that there is a db.relationship between the question, survey and the survey question'

$ws.Range("G41").Value = 'Passed'


# Move the active selection to reflect where editing left off.
$ws.Range("H23").Select()
